$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the now-obsolete trailing rows (old rows 4-8) first, since the
# remaining rows 2 & 3 are simply being updated in place.
$ws.Range("A4:E8").EntireRow.Delete()

# Row 2 updates (A2 stays 1)
$ws.Range("B2").Value = 11
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1000000000000
$ws.Range("E2").Value = 45738.18789351852

# Row 3 updates (A3 stays 2)
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 180
$ws.Range("E3").Value = 45740.53111111111
